$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text cells that look numeric need a leading apostrophe so Excel keeps
# them as text (preserving leading zeros), then reset the style so no
# stray quote-prefix formatting sticks to the cell.
$ws.Range("C2").Value = "'016015"
$ws.Range("C2").Style = "Normal"

$ws.Range("F2").Value = "通讯行业"

$ws.Range("J2").Value = "'001"
$ws.Range("J2").Style = "Normal"

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 141949969.33
$ws.Range("P2").Value = 1372289.09
$ws.Range("Q2").Value = 70656519.59999999
$ws.Range("S2").Value = 35232251.34
$ws.Range("U2").Value = 22857516.57
$ws.Range("W2").Value = 32848460.38
$ws.Range("X2").Value = 17502019.36
$ws.Range("Z2").Value = 4015091.91
$ws.Range("AB2").Value = 109101508.95
$ws.Range("AF2").Value = 422.574935337
$ws.Range("AG2").Value = 23.1408717699
